$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "Player Info" sheet in front of the existing sheets and
#    populate it with the player's bio data.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$infoSheet = $wb.Worksheets.Add($battingSheet)
$infoSheet.Name = "Player Info"

# Sheet handles grabbed before the Add() above can now point at the wrong
# worksheet (handles resolve by position, and Add() shifted everything
# after it along by one) - re-resolve them by name now that the sheet
# collection is stable again.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Header row
$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

$infoHeader = $infoSheet.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160
$infoHeader.Borders.LineStyle = 1

# Data row - ID must stay a text value ("4020"), not a number.
$infoSheet.Range("A2").NumberFormat = "@"
$infoSheet.Range("A2").Value = "4020"
$infoSheet.Range("A2").Style = "Normal"

$infoSheet.Range("B2").Value = "Michael James Gratton Rippon"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Left Arm Wrist Spin (Chinaman)"

$infoSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Rename MATCH_CARD_LINK -> MATCH_CODE and replace the full scorecard URL
#    with just the numeric match code on both the batting and bowling sheets.
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("3526", "3528", "3605", "3610", "4184", "4185", "4563", "4566", "4568")
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $battingSheet.Range("D" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
    $cell.Style = "Normal"
}

$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $bowlingSheet.Range("B" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
    $cell.Style = "Normal"
}
